# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values (column G) recalculated for rows 2-32
$kValues = @(2, 1, 8, 2, 8, 3, 7, 4, 5, 2, 2, 11, 4, 3, 7, 6, 5, 3, 10, 8, 2, 5, 11, 3, 2, 5, 3, 6, 3, 2, 2)

$row = 2
foreach ($val in $kValues) {
    $ws.Range("G$row").Value = $val
    $row = $row + 1
}
